$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting existing rows 32:58 down to 33:59
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly price record
$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value = 44634
$ws.Cells.Item(32, 5).Value = 15
$ws.Cells.Item(32, 6).Value = 100112027
$ws.Cells.Item(32, 7).Value = "Melón"
$ws.Cells.Item(32, 8).Value = "Tuna"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 60
$ws.Cells.Item(32, 11).Value = 16000
$ws.Cells.Item(32, 12).Value = 17000
$ws.Cells.Item(32, 13).Value = 16500
$ws.Cells.Item(32, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(32, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value = 917
$ws.Cells.Item(32, 17).Value = 18
$ws.Cells.Item(32, 18).Value = "Hortaliza"
